# Update the answer key table: replace each three-digit-by-one-digit
# division "problem, remainder" string with the newly generated one.
#
# NOTE: this runtime's Find.Execute operates against the whole document
# story (it is not clipped to the calling Range), so every "old" string
# must be unique in the document at the moment it is searched for. One
# pair of cells collides -- the text produced for row 1 / col 3
# ("566÷7=80, 6") is identical to the *original* text of row 2 / col 2,
# which itself is being replaced later. Doing that row-2/col-2 lookup
# first (while "566÷7=80, 6" still uniquely identifies it) avoids the
# later ReplaceAll accidentally also rewriting the freshly-written
# row-1/col-3 cell.

$d = $word.ActiveDocument

$wdReplaceAll = 2

$d.Content.Find.Execute("305÷4=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "717÷2=358, 1", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("859÷8=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "431÷5=86, 1", $wdReplaceAll) | Out-Null

# (row 2 col 2) must run before (row 1 col 3) -- see note above.
$d.Content.Find.Execute("566÷7=80, 6", $true, $false, $false, $false, $false, $true, 1, $false, "212÷6=35, 2", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("315÷9=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "566÷7=80, 6", $wdReplaceAll) | Out-Null

$d.Content.Find.Execute("249÷4=62, 1", $true, $false, $false, $false, $false, $true, 1, $false, "914÷7=130, 4", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("963÷7=137, 4", $true, $false, $false, $false, $false, $true, 1, $false, "264÷7=37, 5", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("764÷8=95, 4", $true, $false, $false, $false, $false, $true, 1, $false, "874÷7=124, 6", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("809÷2=404, 1", $true, $false, $false, $false, $false, $true, 1, $false, "387÷9=43, 0", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("365÷7=52, 1", $true, $false, $false, $false, $false, $true, 1, $false, "648÷6=108, 0", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("234÷6=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "402÷9=44, 6", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("242÷6=40, 2", $true, $false, $false, $false, $false, $true, 1, $false, "651÷2=325, 1", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("651÷7=93, 0", $true, $false, $false, $false, $false, $true, 1, $false, "995÷4=248, 3", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("240÷4=60, 0", $true, $false, $false, $false, $false, $true, 1, $false, "781÷6=130, 1", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("103÷7=14, 5", $true, $false, $false, $false, $false, $true, 1, $false, "698÷5=139, 3", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("463÷6=77, 1", $true, $false, $false, $false, $false, $true, 1, $false, "157÷4=39, 1", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("695÷5=139, 0", $true, $false, $false, $false, $false, $true, 1, $false, "695÷9=77, 2", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("837÷3=279, 0", $true, $false, $false, $false, $false, $true, 1, $false, "870÷7=124, 2", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("581÷2=290, 1", $true, $false, $false, $false, $false, $true, 1, $false, "461÷2=230, 1", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("340÷7=48, 4", $true, $false, $false, $false, $false, $true, 1, $false, "857÷2=428, 1", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("142÷6=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "639÷6=106, 3", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("869÷4=217, 1", $true, $false, $false, $false, $false, $true, 1, $false, "694÷7=99, 1", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("693÷6=115, 3", $true, $false, $false, $false, $false, $true, 1, $false, "279÷5=55, 4", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("678÷7=96, 6", $true, $false, $false, $false, $false, $true, 1, $false, "540÷5=108, 0", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("210÷5=42, 0", $true, $false, $false, $false, $false, $true, 1, $false, "282÷8=35, 2", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("899÷5=179, 4", $true, $false, $false, $false, $false, $true, 1, $false, "830÷7=118, 4", $wdReplaceAll) | Out-Null
